# Session 5: Dynamic Programming
# - Grade the "Session 5 (Dynamic Prog.)" column (F) for row 4 with a 9.
# - Leave a grading comment in F5 (merged F5:F12), matching the wrapped
#   comment style already used in the other feedback cells of row 5.
# - Move the active selection onto the column that was just filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F4: mark for Session 5 (Dynamic Programming). J4/K4 recalc automatically.
$ws.Range("F4").Value = 9

# F5 (merged F5:F12): grading comment for that session.
$ws.Range("F5").Value = "Very good. Stack memory is around O(1) for dynamic programming and O(n) for the recursive version of the algorithm (taking into account the height of the tree of states). Recursive implementation should be O(3^n) since we need to always do the 3 calls to guarantee a correct solution."

# Match the wrapped, top/left-aligned look used by the other comment cells.
$ws.Range("F5:F12").WrapText = $true
$ws.Range("F5:F12").HorizontalAlignment = -4131   # xlLeft
$ws.Range("F5:F12").VerticalAlignment = -4160     # xlTop

# The author's selection ends up on the column just filled in.
$ws.Range("F5:F12").Select()
